# parameters.xlsx — add Carbon dioxide uptake kinetic parameters,
# fix a couple of parameter values, and tidy two legacy cells.
# (Commit: "Fixed critical bug in f function" / loop changed protein
#  producing parameters but never re-ran FBA afterwards — this edit
#  patches the parameter sheet that feeds that loop.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Clarify the units on the glucose inflow concentration (gluf):
#    "mM" -> "mM/L"  (row 110, column D, before the new rows shift
#    everything down)
# ---------------------------------------------------------------
$ws.Cells.Item(110, 4).Value = "mM/L"

# ---------------------------------------------------------------
# 2. Insert two new rows right before the old row 75 (the start of
#    the "Initial biomass concentrations" block) to hold the new
#    "Carbon dioxide uptake kinetic parameters" section.
#    This naturally inherits the formatting of the row above.
# ---------------------------------------------------------------
$ws.Rows.Item(75).Resize(2).Insert()

$ws.Cells.Item(75, 1).Value = "vmax_c4"
$ws.Cells.Item(76, 1).Value = "ks_c4"

$ws.Cells.Item(75, 2).Value = 2.5
$ws.Cells.Item(76, 2).Value = 0.02

$ws.Cells.Item(75, 3).Value = "%Carbon dioxide uptake kinetic parameters"

# ---------------------------------------------------------------
# 3. Parameter-value corrections (row numbers below already account
#    for the two rows inserted above):
#      glu0 (initial glucose concentration)   0   -> 5
#      sbof (inflow biomass sbo)              1   -> 10
#      oxyf (inflow oxygen concentration)   150   -> 500
#      phof (inflow phosphate concentration)250   -> 500
# ---------------------------------------------------------------
$ws.Cells.Item(83, 2).Value = 5
$ws.Cells.Item(106, 2).Value = 10
$ws.Cells.Item(114, 2).Value = 500
$ws.Cells.Item(115, 2).Value = 500

# ---------------------------------------------------------------
# 4. Row 11 (vmax_o1): an (empty) D cell picks up the same style as
#    the rest of the row.
# ---------------------------------------------------------------
$ws.Range("C11").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Cells.Item(11, 4).ClearContents()

# ---------------------------------------------------------------
# 5. Row 29 (vmax_am1): the value used to live in D29 (unused helper
#    column) - move it into B29 (the actual parameter value) with
#    matching formatting, and drop the now-redundant D29 cell.
# ---------------------------------------------------------------
$ws.Range("D29").Copy()
$ws.Range("B29").PasteSpecial(-4122)
$ws.Cells.Item(29, 2).Value = 0.623
$ws.Cells.Item(29, 4).Clear()

# ---------------------------------------------------------------
# 6. The hidden AutoFilter helper name needs to grow by the same two
#    rows we just inserted (133 -> 135).
# ---------------------------------------------------------------
$name = $wb.Names.Item(1)
$name.RefersTo = "=Sheet1!`$A`$1:`$D`$135"

# ---------------------------------------------------------------
# 7. Leave the view where the edits were made.
# ---------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A77").Select()
